$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    16 = 'JUSTICE CATALYST LAW'
    24 = 'LATHAM & WATKINS'
    25 = 'CLARK GROUP'
    28 = 'CLIFFORD CHANCE US'
    32 = 'WILEY REIN'
    35 = 'CNKINDER'
    36 = 'GLOBAL COUNSEL'
    47 = 'VISA'
    54 = 'MANATT PHELPS AND PHILLIPS'
    57 = 'OTHERSIDE CONSULTING'
    68 = 'IMBUE GROUP'
    76 = 'CACI INTERNATIONAL'
    78 = 'CLOUDFLARE'
    81 = 'SENSHI AME ADVISORS'
    82 = 'INSTITUTE OF ELECTRONIC AND ELECTRICAL ENGINEERS'
    103 = 'TEACH PLUS'
    107 = 'RABEN GROUP'
    124 = 'HIGH WATER MARK'
    128 = 'BFPCC'
    129 = 'MANATT PHELPS AND PHILLIPS'
    133 = 'BOIES SCHILLER FLEXNER'
    141 = 'HOGAN LOVELLS'
    160 = 'INSTITUTES FOR BEHAVIOR RESOURCES'
    163 = 'GUIDEHOUSE'
    175 = 'MANATT PHELPS AND PHILLIPS'
    180 = 'TOFFLER ASSOCIATES'
    183 = 'ONE CONCERN'
    187 = 'SHEPPARD MULLIN RICHTER & HAMPTON'
    193 = 'LEVY GROUP'
    197 = 'CROWELL & MORING'
    200 = 'PAUL HASTINGS'
    202 = 'ADVOCACY BLUEPRINTS'
    227 = 'HARRIS WILTSHIRE AND GRANNIS'
    228 = 'LINKLATERS'
    229 = 'LATHAM & WATKINS'
    231 = 'ARNOLD VENTURES'
    233 = 'WILMERHALE'
    234 = 'MAYER BROWN'
    235 = 'MORGAN LEWIS & BOCKIUS'
    238 = 'COOLEY'
    245 = 'DAVIS POLK & WARDWELL'
    254 = 'HUNTON ANDREWS KURTH'
    280 = 'JOBS FOR THE FUTURE'
    303 = 'HWC'
    313 = 'FACEBOOK'
    331 = 'SKADDEN ARPS SLATE MEAGHER & FLOM'
    335 = 'INTERNATIONAL CAPITAL STRATEGIES'
    345 = 'BRIGHTSIDE BENEFIT'
    354 = 'ARNOLD & PORTER KAYE SCHOLER'
    377 = 'ATLAS RESEARCH'
    388 = 'WORKMERK'
    395 = 'LATHAM & WATKINS'
    396 = 'GRIST MAGAZINE'
    412 = 'DLA PIPER'
    423 = 'BETTER MARKETS'
    428 = 'ORRICK HERRINGTON & SUTCLIFFE'
    432 = 'ELEVATE'
    437 = 'CLIO STRATEGIES'
    450 = 'FACEBOOK'
    451 = 'MORRISON & FOERSTER'
    461 = 'VELEZ SECURITY SOLUTIONS'
    468 = 'WHITE & CASE'
    471 = 'SUMS 15 HOLDINGS'
    481 = 'FACEBOOK'
    495 = 'SIDLEY AUSTIN'
    496 = 'AIRBNB'
    503 = 'SEQUOIA CAPITAL OPERATIONS'
    510 = 'LYFT'
    513 = 'AIRBNB'
    514 = 'MEOW WOLF'
    523 = 'LEVY GROUP'
    533 = 'COVINGTON & BURLING'
    554 = 'VISA'
    564 = 'FUTURE PARTNERS'
    573 = 'FACEBOOK'
    582 = 'JUGGERNAUT'
    597 = 'NATIONAL RESOURCES DEFENSE COUNCIL'
    598 = 'BANDERA ELECTRIC COOPERATIVE'
    612 = 'TRUEPIC'
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item([int]$row, 6).Value = $updates[$row]
}
